# Apply the update described by the commit:
# - Rename "Requested quantity" header to "Weekly_PO_Qty" on the Weekly Quantity sheet
# - Rename "Requested quantity" header to "Monthly_PO_Qty" on the Monthly Trend sheet
# - Add a new "PO Forecast" worksheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Rename header cells -----------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header format (bold + border + centered) from the
# Weekly Quantity sheet's header row, and the existing date number format
# from its date column, so no new styles are introduced.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122) # xlPasteFormats

# --- Header row -----------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows --------------------------------------------------------
$data = @(
    @(45578.99999999999, 160, 160.0006889960679,  160.0006890563822),
    @(45592.99999999999, 16,  16.00068899520499,  16.00068905872806),
    @(45599.99999999999, 0,   -55.99931136720888, -55.99931065957059),
    @(45606.99999999999, 0,   -127.9993122417659, -127.9993098772373),
    @(45613.99999999999, 0,   -199.9993136786671, -199.9993086925135),
    @(45620.99999999999, 0,   -271.999315347539,  -271.9993071267634),
    @(45627.99999999999, 0,   -343.9993168773421, -343.9993054940862),
    @(45634.99999999999, 0,   -415.9993189294993, -415.9993037246624),
    @(45641.99999999999, 0,   -487.9993214387752, -487.9993019985576),
    @(45648.99999999999, 0,   -559.9993238781843, -559.999300072659)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
